$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-9 data, columns A-T (FAPs/sCs sending-cluster split, refreshed LR-stats)
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rspo3"
$ws.Range("C2").Value = "Fzd8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.83236
$ws.Range("H2").Value = 14.49708
$ws.Range("I2").Value = 0.975350813525687
$ws.Range("J2").Value = 0.975350813525687
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.554362333333333
$ws.Range("N2").Value = 10.663087
$ws.Range("O2").Value = 0.2019894022634335
$ws.Range("P2").Value = 0.2019894022634335
$ws.Range("Q2").Value = 17.17595836510667
$ws.Range("R2").Value = 154.58362528596
$ws.Range("S2").Value = 0.1970105278212071
$ws.Range("T2").Value = 0.1970105278212071

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo3"
$ws.Range("C3").Value = "Fzd8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.83236
$ws.Range("H3").Value = 14.49708
$ws.Range("I3").Value = 0.975350813525687
$ws.Range("J3").Value = 0.975350813525687
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.625787
$ws.Range("N3").Value = 34.877361
$ws.Range("O3").Value = 0.6606770910634029
$ws.Range("P3").Value = 0.6606770910634029
$ws.Range("Q3").Value = 56.17998806732001
$ws.Range("R3").Value = 505.61989260588
$ws.Range("S3").Value = 0.6443919382464743
$ws.Range("T3").Value = 0.6443919382464743

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Fzd8"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.83236
$ws.Range("H4").Value = 14.49708
$ws.Range("I4").Value = 0.975350813525687
$ws.Range("J4").Value = 0.975350813525687
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02609533333333333
$ws.Range("N4").Value = 0.078286
$ws.Range("O4").Value = 0.001482961017348462
$ws.Range("P4").Value = 0.001482961017348462
$ws.Range("Q4").Value = 0.1261020449866667
$ws.Range("R4").Value = 1.13491840488
$ws.Range("S4").Value = 0.001446407234697702
$ws.Range("T4").Value = 0.001446407234697702

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo3"
$ws.Range("C5").Value = "Fzd8"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.83236
$ws.Range("H5").Value = 14.49708
$ws.Range("I5").Value = 0.975350813525687
$ws.Range("J5").Value = 0.975350813525687
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.390531666666666
$ws.Range("N5").Value = 7.171595
$ws.Range("O5").Value = 0.1358505456558151
$ws.Range("P5").Value = 0.1358505456558151
$ws.Range("Q5").Value = 11.55190960473333
$ws.Range("R5").Value = 103.9671864426
$ws.Range("S5").Value = 0.1325019402233077
$ws.Range("T5").Value = 0.1325019402233077

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Rspo3"
$ws.Range("C6").Value = "Fzd8"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.122124
$ws.Range("H6").Value = 0.366372
$ws.Range("I6").Value = 0.02464918647431296
$ws.Range("J6").Value = 0.02464918647431296
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.554362333333333
$ws.Range("N6").Value = 10.663087
$ws.Range("O6").Value = 0.2019894022634335
$ws.Range("P6").Value = 0.2019894022634335
$ws.Range("Q6").Value = 0.434072945596
$ws.Range("R6").Value = 3.906656510363999
$ws.Range("S6").Value = 0.004978874442226384
$ws.Range("T6").Value = 0.004978874442226386

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Rspo3"
$ws.Range("C7").Value = "Fzd8"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.122124
$ws.Range("H7").Value = 0.366372
$ws.Range("I7").Value = 0.02464918647431296
$ws.Range("J7").Value = 0.02464918647431296
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 11.625787
$ws.Range("N7").Value = 34.877361
$ws.Range("O7").Value = 0.6606770910634029
$ws.Range("P7").Value = 0.6606770910634029
$ws.Range("Q7").Value = 1.419787611588
$ws.Range("R7").Value = 12.778088504292
$ws.Range("S7").Value = 0.01628515281692846
$ws.Range("T7").Value = 0.01628515281692847

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Rspo3"
$ws.Range("C8").Value = "Fzd8"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.122124
$ws.Range("H8").Value = 0.366372
$ws.Range("I8").Value = 0.02464918647431296
$ws.Range("J8").Value = 0.02464918647431296
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02609533333333333
$ws.Range("N8").Value = 0.078286
$ws.Range("O8").Value = 0.001482961017348462
$ws.Range("P8").Value = 0.001482961017348462
$ws.Range("Q8").Value = 0.003186866488
$ws.Range("R8").Value = 0.028681798392
$ws.Range("S8").Value = 0.00003655378265075908
$ws.Range("T8").Value = 0.00003655378265075909

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Rspo3"
$ws.Range("C9").Value = "Fzd8"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.122124
$ws.Range("H9").Value = 0.366372
$ws.Range("I9").Value = 0.02464918647431296
$ws.Range("J9").Value = 0.02464918647431296
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.390531666666666
$ws.Range("N9").Value = 7.171595
$ws.Range("O9").Value = 0.1358505456558151
$ws.Range("P9").Value = 0.1358505456558151
$ws.Range("Q9").Value = 0.29194128926
$ws.Range("R9").Value = 2.62747160334
$ws.Range("S9").Value = 0.003348605432507353
$ws.Range("T9").Value = 0.003348605432507354
